$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.974.91"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.065.56"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.77"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.62"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.376"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.588.29"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.38"
$ws.Range("E14").Value = "  +4.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000164"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.994.33"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.063.53"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.12"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.50"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.501"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.39"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0908"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.47"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  +6.15%  "
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.62"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.99"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.53"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0679"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.105.78"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.74"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.657"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.299.89"
$ws.Range("E44").Value = "  +4.40%  "
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.86"
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.941"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.93"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.738"
$ws.Range("E50").Value = "  +9.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "253.60"
$ws.Range("E51").Value = "  +9.46%  "
